# RPA datasets push 2024-05-29
# Insert a new IPO demand-forecast record as the new row 2, pushing the
# existing data rows (previously rows 2-13) down to rows 3-14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 2 (shifts rows 2:13 -> 3:14).
$ws.Rows.Item(2).Insert()

# Force the three date columns to be stored as plain text (matching how
# every other date value in this sheet is stored as text, not as an
# Excel date serial number).
$ws.Range("A2:C2").NumberFormat = "@"

# Populate the new row with the latest pushed record.
$ws.Range("A2").Value = "2024-05-13"
$ws.Range("B2").Value = "2024-05-14"
$ws.Range("C2").Value = "2024-05-29"
$ws.Range("D2").Value = "미래"
$ws.Range("E2").Value = "미래에셋비전스팩4호"
$ws.Range("F2").Value = 6650000
$ws.Range("G2").Value = 6650000
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 8100000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 2000
$ws.Range("N2").Value = "1011.2:1"
$ws.Range("O2").Value = "-"
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = "기업인수목적회사(기타금융서비스)"

# The row insert operation picks up incidental formatting (e.g. a bold
# font/border inherited from neighboring cells); reset the new row back
# to the plain/default style used by the rest of the data rows.
$ws.Range("A2:Y2").Style = "Normal"
